$d = $word.ActiveDocument

$pairs = @(
    @("509÷7=72, 5", "428÷3=142, 2"),
    @("101÷8=12, 5", "956÷5=191, 1"),
    @("756÷8=94, 4", "314÷7=44, 6"),
    @("504÷4=126, 0", "922÷5=184, 2"),
    @("111÷9=12, 3", "271÷9=30, 1"),
    @("601÷7=85, 6", "787÷8=98, 3"),
    @("389÷2=194, 1", "564÷6=94, 0"),
    @("997÷3=332, 1", "153÷9=17, 0"),
    @("335÷8=41, 7", "576÷3=192, 0"),
    @("198÷5=39, 3", "903÷6=150, 3"),
    @("985÷8=123, 1", "480÷7=68, 4"),
    @("185÷6=30, 5", "630÷7=90, 0"),
    @("751÷6=125, 1", "733÷5=146, 3"),
    @("350÷6=58, 2", "118÷8=14, 6"),
    @("455÷6=75, 5", "923÷7=131, 6"),
    @("609÷8=76, 1", "213÷7=30, 3"),
    @("365÷3=121, 2", "366÷6=61, 0"),
    @("427÷6=71, 1", "416÷9=46, 2"),
    @("871÷3=290, 1", "858÷6=143, 0"),
    @("173÷7=24, 5", "750÷7=107, 1"),
    @("292÷3=97, 1", "669÷9=74, 3"),
    @("502÷3=167, 1", "501÷7=71, 4"),
    @("102÷4=25, 2", "856÷9=95, 1"),
    @("732÷6=122, 0", "109÷8=13, 5"),
    @("392÷6=65, 2", "129÷2=64, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
